$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 first (SC 92), then row 26 (RM 232), so earlier deletion
# doesn't shift the row index of the later one.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()
